# Helper: find a shape on a slide by its Name property.
function Find-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. Nudge the six "pebble" ellipses horizontally (y unchanged).
# ---------------------------------------------------------------------------
(Find-ShapeByName $s "Ellipse 26").Left = 292.0470275878906
(Find-ShapeByName $s "Ellipse 27").Left = 292.0470275878906
(Find-ShapeByName $s "Ellipse 40").Left = 189.7174835205078
(Find-ShapeByName $s "Ellipse 41").Left = 87.3880386352539
(Find-ShapeByName $s "Ellipse 42").Left = 189.7174835205078
(Find-ShapeByName $s "Ellipse 43").Left = 87.3880386352539

# ---------------------------------------------------------------------------
# 2. Widen / reposition "Rectangle 17" (the left sand block) and the two
#    connectors glued to it.
# ---------------------------------------------------------------------------
$rect17 = Find-ShapeByName $s "Rectangle 17"
$rect17.Left = 62.490394592285156
$rect17.Top = 215.75433349609375
$rect17.Width = 316.4704895019531

$cxn19 = Find-ShapeByName $s "Connecteur droit 19"
$cxn19.Left = 170.22560119628906
$cxn19.Top = 215.75433349609375

$cxn23 = Find-ShapeByName $s "Connecteur droit 23"
$cxn23.Left = 62.490394592285156
$cxn23.Top = 290.0647277832031
$cxn23.Width = 316.4704895019531

# ---------------------------------------------------------------------------
# 3. Drop the two vertical filler rectangles that used to sit mid-tank.
# ---------------------------------------------------------------------------
(Find-ShapeByName $s "Rectangle 22").Delete()
(Find-ShapeByName $s "Rectangle 25").Delete()

# ---------------------------------------------------------------------------
# 4. Nudge the four creature name labels that sit above the ellipses.
# ---------------------------------------------------------------------------
(Find-ShapeByName $s "ZoneTexte 32").Left = 283.8028564453125
(Find-ShapeByName $s "ZoneTexte 33").Left = 190.1252899169922
(Find-ShapeByName $s "ZoneTexte 34").Left = 193.6593780517578
(Find-ShapeByName $s "ZoneTexte 29").Left = 283.6816711425781

# ---------------------------------------------------------------------------
# 5. Add the new vertical divider connector ("Connecteur droit 31", id 32).
#    The engine hands out shape ids from an internal monotonic counter that
#    was seeded from the ids already present when the deck was loaded; we
#    burn through the intermediate ids with throw-away lines (deleted right
#    away) so that the *real* new connector (made via Duplicate so it
#    inherits the exact line/style formatting already used by its sibling
#    connectors) lands on id 32, matching the target deck.
# ---------------------------------------------------------------------------
$throwaways = @()
for ($i = 0; $i -lt 16; $i++) {
    $throwaways += $s.Shapes.AddLine(1, 1, 1, 2)
}
foreach ($t in $throwaways) {
    $t.Delete()
}

$template = Find-ShapeByName $s "Connecteur droit 19"
$newCxnRange = $template.Duplicate()
$newCxn = $newCxnRange.Item(1)
$newCxn.Name = "Connecteur droit 31"
$newCxn.Left = 273.78033447265625
$newCxn.Top = 215.58543395996094
